# ChatTrinityEnglishSourceData.xlsx update:
#  - insert a new "Web" data-source row (FAQ page about registered exam
#    centres) just above the existing "PDF" block, pushing the PDF rows
#    (old rows 136-192) down to 137-193
#  - hyperlink the new row's URL cell (gets the built-in "Hyperlink" style)
#  - grow the AutoFilter / hidden _FilterDatabase range to match the extra row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the PDF rows down one slot by inserting a fresh row at 136.
$ws.Rows.Item(136).Insert()

# Populate the newly inserted row with the Web/FAQ entry.
$ws.Range("A136").Value = "Web"
$ws.Hyperlinks.Add($ws.Range("B136"), "https://www.trinitycollege.com/about-us/work-with-trinity/registered-exam-centre")

# Re-establish the autofilter over the now-larger A1:B193 range.
$ws.AutoFilterMode = $false
$ws.Range("A1:B193").AutoFilter()

# Extend the hidden _FilterDatabase defined name to match.
$wb.Names.Item(1).RefersTo = "=Sheet1!`$A`$1:`$B`$193"

# Keep the active selection tidy (back to the top-left cell).
$ws.Range("A1").Select()
